# Update market price / profit figures across multiple sheets
# as pulled by the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1467.5
$ws.Range("I15").Value = 1467.5
$ws.Range("K15").Value = 4402.5
$ws.Range("M15").Value = -4233.5
$ws.Range("H18").Value = 912.5
$ws.Range("I18").Value = 912.5
$ws.Range("K18").Value = 912.5
$ws.Range("M18").Value = -628.5
$ws.Range("H111").Value = 638
$ws.Range("I111").Value = 561.625
$ws.Range("K111").Value = 1684.875
$ws.Range("M111").Value = 1382.125
$ws.Range("H125").Value = 6836.5
$ws.Range("I125").Value = 4340.3335
$ws.Range("K125").Value = 39063.0015
$ws.Range("M125").Value = -36603.0015
$ws.Range("H138").Value = 3875.1365
$ws.Range("J138").Value = 4259.222
$ws.Range("L138").Value = 12777.666
$ws.Range("N138").Value = -23057.666
$ws.Range("H140").Value = 50780
$ws.Range("J140").Value = 50780
$ws.Range("L140").Value = 50780
$ws.Range("N140").Value = -61140
$ws.Range("H141").Value = 2426.2666
$ws.Range("I141").Value = 2426.2666
$ws.Range("K141").Value = 7278.7998
$ws.Range("M141").Value = -2098.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2147.3333
$ws.Range("I32").Value = 2147.3333
$ws.Range("K32").Value = 2147.3333
$ws.Range("M32").Value = -1860.3333
$ws.Range("H74").Value = 3330.1667
$ws.Range("I74").Value = 3330.1667
$ws.Range("K74").Value = 3330.1667
$ws.Range("M74").Value = -2456.1667
$ws.Range("H77").Value = 3330.1667
$ws.Range("I77").Value = 3330.1667
$ws.Range("K77").Value = 16650.8335
$ws.Range("M77").Value = -12282.8335
$ws.Range("H132").Value = 1928.2142
$ws.Range("I132").Value = 1899.6
$ws.Range("K132").Value = 5698.799999999999
$ws.Range("M132").Value = -3168.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 125295
$ws.Range("J68").Value = 125295
$ws.Range("L68").Value = 125295
$ws.Range("N68").Value = -126917
$ws.Range("H71").Value = 125295
$ws.Range("J71").Value = 125295
$ws.Range("L71").Value = 375885
$ws.Range("N71").Value = -383997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1770.4286
$ws.Range("I16").Value = 1770.4286
$ws.Range("K16").Value = 1770.4286
$ws.Range("M16").Value = -1483.4286
$ws.Range("H17").Value = 13989.25
$ws.Range("J17").Value = 13989.25
$ws.Range("L17").Value = 13989.25
$ws.Range("N17").Value = -14337.25
$ws.Range("H25").Value = 1760
$ws.Range("I25").Value = 1013.3333
$ws.Range("J25").Value = 4000
$ws.Range("K25").Value = 1013.3333
$ws.Range("L25").Value = 4000
$ws.Range("M25").Value = -839.3333
$ws.Range("N25").Value = -4348
$ws.Range("H31").Value = 2074.925
$ws.Range("I31").Value = 1504.9667
$ws.Range("J31").Value = 3784.8
$ws.Range("K31").Value = 1504.9667
$ws.Range("L31").Value = 3784.8
$ws.Range("M31").Value = -1209.9667
$ws.Range("N31").Value = -4374.8
$ws.Range("H34").Value = 2074.925
$ws.Range("I34").Value = 1504.9667
$ws.Range("J34").Value = 3784.8
$ws.Range("K34").Value = 1504.9667
$ws.Range("L34").Value = 3784.8
$ws.Range("M34").Value = -1302.9667
$ws.Range("N34").Value = -4188.8
$ws.Range("H41").Value = 18750
$ws.Range("J41").Value = 20000
$ws.Range("L41").Value = 20000
$ws.Range("N41").Value = -20856
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26472
$ws.Range("H61").Value = 25000
$ws.Range("J61").Value = 25000
$ws.Range("L61").Value = 25000
$ws.Range("N61").Value = -25696
$ws.Range("H113").Value = 1770.4286
$ws.Range("I113").Value = 1770.4286
$ws.Range("K113").Value = 1770.4286
$ws.Range("M113").Value = 399.5714
$ws.Range("H132").Value = 2250.1
$ws.Range("I132").Value = 2111.2222
$ws.Range("K132").Value = 6333.6666
$ws.Range("M132").Value = -3803.6666
$ws.Range("H134").Value = 102542.4
$ws.Range("I134").Value = 126678
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 380034
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -377499
$ws.Range("N134").Value = -23070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 200645.3
$ws.Range("I11").Value = 500755
$ws.Range("J11").Value = 572.1667
$ws.Range("K11").Value = 1502265
$ws.Range("L11").Value = 1716.5001
$ws.Range("M11").Value = -1502125
$ws.Range("N11").Value = -1996.5001
$ws.Range("H12").Value = 59.705883
$ws.Range("J12").Value = 86.545456
$ws.Range("L12").Value = 259.636368
$ws.Range("N12").Value = -605.6363679999999
$ws.Range("H63").Value = 1166
$ws.Range("I63").Value = 999.5
$ws.Range("J63").Value = 1499
$ws.Range("K63").Value = 2998.5
$ws.Range("L63").Value = 4497
$ws.Range("M63").Value = -2249.5
$ws.Range("N63").Value = -5995
$ws.Range("H66").Value = 1166
$ws.Range("I66").Value = 999.5
$ws.Range("J66").Value = 1499
$ws.Range("K66").Value = 8995.5
$ws.Range("L66").Value = 13491
$ws.Range("M66").Value = -5251.5
$ws.Range("N66").Value = -20979
$ws.Range("H80").Value = 11834.846
$ws.Range("J80").Value = 11749
$ws.Range("L80").Value = 35247
$ws.Range("N80").Value = -37119
$ws.Range("H83").Value = 11834.846
$ws.Range("J83").Value = 11749
$ws.Range("L83").Value = 105741
$ws.Range("N83").Value = -115101
$ws.Range("H139").Value = 3721.5715
$ws.Range("I139").Value = 1164.7142
$ws.Range("K139").Value = 3494.1426
$ws.Range("M139").Value = 1645.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3997.25
$ws.Range("J70").Value = 3997.25
$ws.Range("L70").Value = 3997.25
$ws.Range("N70").Value = -4537.25
$ws.Range("H73").Value = 3997.25
$ws.Range("J73").Value = 3997.25
$ws.Range("L73").Value = 3997.25
$ws.Range("N73").Value = -5869.25
$ws.Range("H95").Value = 12672
$ws.Range("J95").Value = 12672
$ws.Range("L95").Value = 12672
$ws.Range("N95").Value = -18164
$ws.Range("H97").Value = 972.5
$ws.Range("I97").Value = 972.5
$ws.Range("K97").Value = 972.5
$ws.Range("M97").Value = -476.5
$ws.Range("H132").Value = 2788.1
$ws.Range("I132").Value = 2874.111
$ws.Range("K132").Value = 8622.332999999999
$ws.Range("M132").Value = -6092.332999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 919.3333
$ws.Range("I22").Value = 1014.1
$ws.Range("J22").Value = 729.8
$ws.Range("K22").Value = 1014.1
$ws.Range("L22").Value = 729.8
$ws.Range("M22").Value = -719.1
$ws.Range("N22").Value = -1319.8
$ws.Range("H27").Value = 919.3333
$ws.Range("I27").Value = 1014.1
$ws.Range("J27").Value = 729.8
$ws.Range("K27").Value = 1014.1
$ws.Range("L27").Value = 729.8
$ws.Range("M27").Value = -907.1
$ws.Range("N27").Value = -943.8
$ws.Range("H55").Value = 240.44444
$ws.Range("I55").Value = 158.45454
$ws.Range("K55").Value = 158.45454
$ws.Range("M55").Value = 14.54545999999999
$ws.Range("H101").Value = 22172
$ws.Range("J101").Value = 22172
$ws.Range("L101").Value = 22172
$ws.Range("N101").Value = -28662
$ws.Range("H132").Value = 6553.0713
$ws.Range("J132").Value = 7481.778
$ws.Range("L132").Value = 22445.334
$ws.Range("N132").Value = -27505.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5711.615
$ws.Range("J62").Value = 5909.1816
$ws.Range("L62").Value = 5909.1816
$ws.Range("N62").Value = -7157.1816
$ws.Range("H65").Value = 5711.615
$ws.Range("J65").Value = 5909.1816
$ws.Range("L65").Value = 29545.908
$ws.Range("N65").Value = -35785.908
